# Update "想去人数" (interested-count) values in the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("展览")
$sheet1.Range("F3").Value  = 201
$sheet1.Range("F5").Value  = 5212
$sheet1.Range("F9").Value  = 584
$sheet1.Range("F10").Value = 535
$sheet1.Range("F13").Value = 1439
$sheet1.Range("F14").Value = 4154
$sheet1.Range("F15").Value = 430
$sheet1.Range("F16").Value = 171
$sheet1.Range("F17").Value = 149
$sheet1.Range("F18").Value = 93
$sheet1.Range("F19").Value = 3193
$sheet1.Range("F20").Value = 158
$sheet1.Range("F21").Value = 1064
$sheet1.Range("F24").Value = 188
$sheet1.Range("F25").Value = 99
$sheet1.Range("F29").Value = 289
$sheet1.Range("F30").Value = 25
$sheet1.Range("F31").Value = 53
$sheet1.Range("F33").Value = 13
$sheet1.Range("F34").Value = 13

$sheet4 = $wb.Worksheets.Item("全部类型")
$sheet4.Range("F3").Value  = 201
$sheet4.Range("F6").Value  = 5212
$sheet4.Range("F10").Value = 584
$sheet4.Range("F11").Value = 535
$sheet4.Range("F14").Value = 1439
$sheet4.Range("F15").Value = 4154
$sheet4.Range("F16").Value = 430
$sheet4.Range("F17").Value = 171
$sheet4.Range("F18").Value = 149
$sheet4.Range("F19").Value = 93
$sheet4.Range("F20").Value = 3193
$sheet4.Range("F21").Value = 158
$sheet4.Range("F22").Value = 1064
$sheet4.Range("F25").Value = 188
$sheet4.Range("F26").Value = 99
$sheet4.Range("F30").Value = 289
$sheet4.Range("F31").Value = 25
$sheet4.Range("F32").Value = 53
$sheet4.Range("F34").Value = 13
$sheet4.Range("F35").Value = 13

$wb.Save()
